# Trade #114 closed at 2026-02-17 16:03:08 - unknown UNKNOWN +0.000%
# Updates Summary + Strategy Status aggregate numbers and appends the new
# closed trade row to both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1198.84     # Current Capital
$wsSummary.Range("B4").Value = -1.17       # Total P&L $
$wsSummary.Range("B6").Value = 114         # Total Trades
$wsSummary.Range("B8").Value = 58          # Losing Trades
$wsSummary.Range("B9").Value = 35.09       # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet (MarketMaking row = row 4)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 98.84        # Capital
$wsStatus.Range("D4").Value = 114          # Trades
$wsStatus.Range("E4").Value = -1.17        # P&L $
$wsStatus.Range("F4").Value = -1.16        # P&L %
$wsStatus.Range("G4").Value = 35.09        # Win Rate %

# ---------------------------------------------------------------------
# 3. Append new trade row (115) to "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A115").Value = 114

    # Date / Time columns are stored as plain text in this workbook.
    # Prefix with a leading apostrophe so Excel keeps the look-alike
    # strings as text instead of converting them to date/time serials
    # (this mirrors how the existing rows already store these values).
    $ws.Range("B115").Value = "'2026-02-17"
    $ws.Range("C115").Value = "'16:03:01"

    $ws.Range("D115").Value = "MarketMaking"
    $ws.Range("E115").Value = "DOWN"
    $ws.Range("F115").Value = 0.17
    $ws.Range("G115").Value = 0.16
    $ws.Range("H115").Value = "CLOSED"
    $ws.Range("I115").Value = -5.8824
    $ws.Range("J115").Value = -0.01
    $ws.Range("K115").Value = 98.84
    $ws.Range("L115").Value = 0
    $ws.Range("M115").Value = 0
    $ws.Range("N115").Value = 0.6
    $ws.Range("O115").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P115").Value = "early_exit"
    $ws.Range("Q115").Value = 0.15
}
